# Backup QR Scanner data - append the latest scan/manual log entry as a new
# row at the bottom of the "Neurology" log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row goes right after the current last row of data (row 85 -> row 86).
$newRow = 86

# Column A ("Student ID") holds a purely-numeric-looking value. A plain
# `.Value = "201560"` assignment would be auto-coerced to a Number by the
# engine, but the sheet stores every column as literal text (matches the
# rest of the log). Entering it as a quote-prefixed literal (like typing
# '201560 into Excel) forces a genuine Text cell instead, then resetting the
# cell style back to Normal drops the quote-prefix formatting so no visual
# style changes stick around.
$ws.Cells.Item($newRow, 1).Formula = "'201560"
$ws.Cells.Item($newRow, 1).Style = "Normal"

# The remaining columns are not purely numeric, so plain value assignment
# already keeps them as text (matching "Subject"/"Type"/"User" elsewhere).
$ws.Cells.Item($newRow, 2).Value = "Neurology"
$ws.Cells.Item($newRow, 3).Value = "29/12/2025"
$ws.Cells.Item($newRow, 4).Value = "13:57:00"
$ws.Cells.Item($newRow, 5).Value = "Manual"
$ws.Cells.Item($newRow, 6).Value = "emp17.farah.a.youssef@gmail.com"
